# Misc changes to existing cal sheets:
#  - Changed references to GI05MOAS to GA05MOAS (Ref Des column, sheet Asset_Cal_Info)
#  - Inserted missing Cal parameter name correction for PG562/PG563 row (cleared stale
#    "No calibration coefficient" note in I6 now that a real coefficient value is present)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Asset_Cal_Info")

# Rows in column A that still carry the old "GI05MOAS-..." Ref Des text.
# Processed top-to-bottom so new shared-string entries are appended in the
# same order they first appear in the sheet.
$fixups = @(
    @{ Row = 2;  Value = "GA05MOAS-PG563-01-CTDGVM000" },
    @{ Row = 4;  Value = "GA05MOAS-PG563-02-DOSTAM000" },
    @{ Row = 6;  Value = "GA05MOAS-PG563-05-NUTNRM000" },
    @{ Row = 7;  Value = "GA05MOAS-PG563-05-NUTNRM000" },
    @{ Row = 8;  Value = "GA05MOAS-PG563-05-NUTNRM000" },
    @{ Row = 9;  Value = "GA05MOAS-PG563-05-NUTNRM000" },
    @{ Row = 10; Value = "GA05MOAS-PG563-05-NUTNRM000" },
    @{ Row = 11; Value = "GA05MOAS-PG563-05-NUTNRM000" },
    @{ Row = 12; Value = "GA05MOAS-PG563-05-NUTNRM000" },
    @{ Row = 14; Value = "GA05MOAS-PG563-06-PARADM000" },
    @{ Row = 16; Value = "GA05MOAS-PG563-03-FLORTM000" },
    @{ Row = 17; Value = "GA05MOAS-PG563-03-FLORTM000" },
    @{ Row = 18; Value = "GA05MOAS-PG563-03-FLORTM000" },
    @{ Row = 19; Value = "GA05MOAS-PG563-03-FLORTM000" },
    @{ Row = 21; Value = "GA05MOAS-PG563-04-FLORTO000" },
    @{ Row = 22; Value = "GA05MOAS-PG563-04-FLORTO000" },
    @{ Row = 23; Value = "GA05MOAS-PG563-04-FLORTO000" },
    @{ Row = 24; Value = "GA05MOAS-PG563-04-FLORTO000" },
    @{ Row = 26; Value = "GA05MOAS-PG563-00-ENG000000" }
)

foreach ($fix in $fixups) {
    $ws.Cells.Item($fix.Row, 1).Value = $fix.Value
}

# Clear the stale "No calibration coefficient" note in I6 - row 6 (PG562/CTDGVM
# calibration) actually has a real coefficient value (20.11) in H6, so the note
# was incorrect and is removed.
$ws.Range("I6").Value = ""

# Restore the last active selection on the sheet.
[void]$ws.Range("E31").Select()
